# Auto-generated Excel COM-interop script
# Commit: Update automàtic: dades i banners [2026-02-28 05:20]
# Applies per-cell text updates to the Dades_Meteo sheet matching the
# upstream meteo.cat re-extraction diff (refreshed DATA_EXTRACCIO timestamps
# plus the handful of measurement values that shifted on re-scrape).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dades_Meteo")

$ws.Range('E2').Value = '2026-02-28 05:18:16'
$ws.Range('N2').Value = '0.1 °C 4:41 TU'
$ws.Range('O2').Value = '0.9 °C'
$ws.Range('E3').Value = '2026-02-28 05:18:18'
$ws.Range('O3').Value = '-0.3 °C'
$ws.Range('E4').Value = '2026-02-28 05:18:20'
$ws.Range('O4').Value = '7.8 °C'
$ws.Range('E5').Value = '2026-02-28 05:18:23'
$ws.Range('N5').Value = '-0.8 °C 4:56 TU'
$ws.Range('O5').Value = '-0.1 °C'
$ws.Range('E6').Value = '2026-02-28 05:18:25'
$ws.Range('E7').Value = '2026-02-28 05:18:27'
$ws.Range('H7').Value = '''77%'
$ws.Range('O7').Value = '12.3 °C'
$ws.Range('E8').Value = '2026-02-28 05:18:30'
$ws.Range('J8').Value = '1023.8 hPa'
$ws.Range('M8').Value = '8.7 °C 4:55 TU'
$ws.Range('E9').Value = '2026-02-28 05:18:32'
$ws.Range('E10').Value = '2026-02-28 05:18:34'
$ws.Range('M10').Value = '8.6 °C 4:46 TU'
$ws.Range('O10').Value = '7.7 °C'
$ws.Range('E11').Value = '2026-02-28 05:18:37'
$ws.Range('N11').Value = '2.1 °C 4:50 TU'
$ws.Range('O11').Value = '3.5 °C'
$ws.Range('E12').Value = '2026-02-28 05:18:38'
$ws.Range('E13').Value = '2026-02-28 05:18:41'
$ws.Range('H13').Value = '''87%'
$ws.Range('J13').Value = '1026.2 hPa'
$ws.Range('N13').Value = '-0.9 °C 4:41 TU'
$ws.Range('O13').Value = '1.3 °C'
$ws.Range('E14').Value = '2026-02-28 05:18:43'
$ws.Range('M14').Value = '11.5 °C 4:55 TU'
$ws.Range('O14').Value = '10.2 °C'
$ws.Range('E15').Value = '2026-02-28 05:18:46'
$ws.Range('E16').Value = '2026-02-28 05:18:48'
$ws.Range('H16').Value = '''61%'
$ws.Range('N16').Value = '-1.5 °C 4:36 TU'
$ws.Range('E17').Value = '2026-02-28 05:18:50'
$ws.Range('G17').Value = '1 cm'
$ws.Range('N17').Value = '3.8 °C 4:40 TU'
$ws.Range('E18').Value = '2026-02-28 05:18:53'
$ws.Range('M18').Value = '9.5 °C 4:56 TU'
$ws.Range('O18').Value = '8.2 °C'
$ws.Range('E19').Value = '2026-02-28 05:18:55'
$ws.Range('H19').Value = '''65%'
$ws.Range('L19').Value = '25.6 km/h - 68º 4:43 TU'
$ws.Range('E20').Value = '2026-02-28 05:18:57'
$ws.Range('N20').Value = '-1.2 °C 4:56 TU'
$ws.Range('O20').Value = '0.1 °C'
$ws.Range('E21').Value = '2026-02-28 05:19:00'
$ws.Range('H21').Value = '''76%'
$ws.Range('N21').Value = '3.5 °C 4:46 TU'
$ws.Range('O21').Value = '5.3 °C'
$ws.Range('E22').Value = '2026-02-28 05:19:02'
$ws.Range('H22').Value = '''59%'
$ws.Range('E23').Value = '2026-02-28 05:19:04'
$ws.Range('N23').Value = '-0.9 °C 4:41 TU'
$ws.Range('O23').Value = '0.0 °C'
$ws.Range('E24').Value = '2026-02-28 05:19:06'
$ws.Range('J24').Value = '1023.5 hPa'
$ws.Range('E25').Value = '2026-02-28 05:19:09'
$ws.Range('H25').Value = '''56%'
$ws.Range('L25').Value = '9.0 km/h - 343º 4:50 TU'
$ws.Range('N25').Value = '-0.7 °C 4:43 TU'
$ws.Range('O25').Value = '1.0 °C'
$ws.Range('E26').Value = '2026-02-28 05:19:11'
$ws.Range('H26').Value = '''70%'
$ws.Range('N26').Value = '4.0 °C 4:45 TU'
$ws.Range('E27').Value = '2026-02-28 05:19:14'
$ws.Range('H27').Value = '''36%'
$ws.Range('N27').Value = '0.7 °C 4:56 TU'
$ws.Range('O27').Value = '2.4 °C'
$ws.Range('E28').Value = '2026-02-28 05:19:16'
$ws.Range('J28').Value = '1024.6 hPa'
$ws.Range('E29').Value = '2026-02-28 05:19:18'
$ws.Range('L29').Value = '12.6 km/h - 354º 4:33 TU'
$ws.Range('O29').Value = '8.8 °C'
$ws.Range('E30').Value = '2026-02-28 05:19:21'
$ws.Range('E31').Value = '2026-02-28 05:19:23'
$ws.Range('L31').Value = '41.8 km/h - 5º 4:55 TU'
$ws.Range('N31').Value = '9.6 °C 4:43 TU'
$ws.Range('E32').Value = '2026-02-28 05:19:25'
$ws.Range('H32').Value = '''90%'
$ws.Range('O32').Value = '5.1 °C'
$ws.Range('E33').Value = '2026-02-28 05:19:28'
$ws.Range('H33').Value = '''71%'
$ws.Range('J33').Value = '1023.6 hPa'
$ws.Range('N33').Value = '3.6 °C 4:38 TU'
$ws.Range('O33').Value = '5.1 °C'
$ws.Range('E34').Value = '2026-02-28 05:19:30'
$ws.Range('H34').Value = '''68%'
$ws.Range('N34').Value = '-0.8 °C 4:59 TU'
$ws.Range('E35').Value = '2026-02-28 05:19:33'
$ws.Range('H35').Value = '''82%'
$ws.Range('J35').Value = '1022.8 hPa'
$ws.Range('L35').Value = '28.8 km/h - 274º 4:59 TU'
$ws.Range('N35').Value = '5.5 °C 4:57 TU'
$ws.Range('O35').Value = '6.8 °C'
$ws.Range('E36').Value = '2026-02-28 05:19:35'
$ws.Range('L36').Value = '18.7 km/h - 14º 4:40 TU'
$ws.Range('M36').Value = '11.5 °C 4:48 TU'
$ws.Range('O36').Value = '10.0 °C'
$ws.Range('E37').Value = '2026-02-28 05:19:37'
$ws.Range('N37').Value = '4.0 °C 4:59 TU'
$ws.Range('E38').Value = '2026-02-28 05:19:39'
$ws.Range('O38').Value = '9.0 °C'
$ws.Range('E39').Value = '2026-02-28 05:19:42'
$ws.Range('H39').Value = '''50%'
$ws.Range('N39').Value = '-1.1 °C 4:46 TU'
$ws.Range('O39').Value = '0.2 °C'
$ws.Range('E40').Value = '2026-02-28 05:19:44'
$ws.Range('G40').Value = '3 cm'
$ws.Range('H40').Value = '''93%'
$ws.Range('N40').Value = '2.0 °C 4:52 TU'
$ws.Range('O40').Value = '3.5 °C'
$ws.Range('E41').Value = '2026-02-28 05:19:46'
$ws.Range('E42').Value = '2026-02-28 05:19:48'
$ws.Range('E43').Value = '2026-02-28 05:19:51'
$ws.Range('H43').Value = '''85%'
$ws.Range('N43').Value = '3.0 °C 4:57 TU'
$ws.Range('E44').Value = '2026-02-28 05:19:53'
$ws.Range('L44').Value = '19.8 km/h - 2º 4:30 TU'
$ws.Range('E45').Value = '2026-02-28 05:19:55'
$ws.Range('H45').Value = '''90%'
$ws.Range('L45').Value = '17.6 km/h - 21º 4:44 TU'
$ws.Range('N45').Value = '5.3 °C 4:59 TU'
$ws.Range('O45').Value = '7.0 °C'
$ws.Range('E46').Value = '2026-02-28 05:19:58'
$ws.Range('J46').Value = '1023.1 hPa'
$ws.Range('M46').Value = '11.4 °C 4:41 TU'
